$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7988880276679993
$ws.Range("B1").Value = 1.87778627872467
$ws.Range("C1").Value = 5.037609100341797
$ws.Range("D1").Value = 2.605562925338745
$ws.Range("E1").Value = 1.404710173606873
